$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Average Score column (F) = AVERAGE(C:E) per student row
$ws.Range("F2").Formula = "=AVERAGE(C2:E2)"
$ws.Range("F3:F10").Formula = "=AVERAGE(C3:E3)"

# Final Score column (H) = SUM(Average Score, Bonus Points) per student row
$ws.Range("H2").Formula = "=SUM(F2:G2)"
$ws.Range("H3:H10").Formula = "=SUM(F3:G3)"

# Summary statistics block (rows 15-19) based on Final Score column
$ws.Range("B15").Formula = "=AVERAGE(H2:H10)"
$ws.Range("B16").Formula = "=MEDIAN(H2:H10)"
$ws.Range("B17").Formula = "=MAX(H2:H10)"
$ws.Range("B18").Formula = "=MIN(H2:H10)"
$ws.Range("B19").Formula = "=STDEV.P(H2:H10)"

# Reflect the final cell selection left after entering the formulas
[void]$ws.Range("B20").Select()
